$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells are written in the same order the shared-string table records them,
# so the rebuilt workbook's sharedStrings.xml ends up byte-for-byte aligned
# with the authored one.
$ws.Range("A1").Value = "OneAgainstMany"
$ws.Range("A2").Value = "Name:"
$ws.Range("C2").Value = "Raw Value"
$ws.Range("B2").Value = "Calculated Value"
$ws.Range("D2").Value = "Percentage Of"
$ws.Range("D1").Value = "Value Used:"
$ws.Range("E1").Value = "Here"
$ws.Range("E2").Value = "Index Pos"
$ws.Range("F2").Value = "Index Length"
$ws.Range("G2").Value = "Against Average"
$ws.Range("H2").Value = "Against Median"
$ws.Range("I2").Value = "Standard Deviations Away"
$ws.Range("J2").Value = "Value Rarity"
$ws.Range("K2").Value = "Relevancy"
$ws.Range("L2").Value = "Relevancy * Against {Value}"
$ws.Range("N2").Value = "Stat Mean"
$ws.Range("O2").Value = "Stat Quartile 1"
$ws.Range("P2").Value = "Stat Quartile 2"
$ws.Range("Q2").Value = "Stat Quartile 3"
$ws.Range("R2").Value = "Stat Mode"
$ws.Range("S2").Value = "Stat Standard Deviation"
$ws.Range("T2").Value = "Stat Grouped Mode"

# Column widths for E:F (authored width is 8.88671875 characters)
$ws.Range("E1:F1").ColumnWidth = 8

# Cursor / selection ends on Z2 as in the saved file
$ws.Range("Z2").Select()
